$wb = $excel.ActiveWorkbook

# "Metadata" sheet holds Property/Value pairs; row 4 is "Name" and its
# value cell (B4) was empty, now it should contain "LangueVs".
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B4").Value = "LangueVs"

# Update the recorded generation date/time string (row 8, "Date").
$meta.Range("B8").Value = "2025-07-18T06:40:38+00:00"
